$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: A=111742138 (Blek stjärnmossa)
$ws.Range("A6").Value2 = 111742138
$ws.Range("B6").Value2 = 92683
$ws.Range("D6").Value2 = 'LC'
$ws.Range("E6").Value2 = 2362
$ws.Range("F6").Value2 = 'Blek stjärnmossa'
$ws.Range("G6").Value2 = 'Mnium stellare'
$ws.Range("H6").Value2 = 'Hedw.'
$ws.Range("Q6").Value2 = 331789
$ws.Range("R6").Value2 = 6626790
$ws.Range("L6").Value2 = "'"
$ws.Range("AC6").ClearContents()

# Row 7: A=111742101 (Vedtrappmossa)
$ws.Range("A7").Value2 = 111742101
$ws.Range("B7").Value2 = 94134
$ws.Range("D7").Value2 = 'NT'
$ws.Range("E7").Value2 = 53
$ws.Range("F7").Value2 = 'Vedtrappmossa'
$ws.Range("G7").Value2 = 'Crossocalyx hellerianus'
$ws.Range("H7").Value2 = '(Nees ex Lindenb.) Meyl.'
$ws.Range("Q7").Value2 = 331780
$ws.Range("R7").Value2 = 6626798
$ws.Range("L7").Value2 = "'"
$ws.Range("AC7").ClearContents()

# Row 8: A=111742170 (Vedticka)
$ws.Range("A8").Value2 = 111742170
$ws.Range("B8").Value2 = 89369
$ws.Range("D8").Value2 = 'LC'
$ws.Range("E8").Value2 = 5447
$ws.Range("F8").Value2 = 'Vedticka'
$ws.Range("G8").Value2 = 'Fuscoporia viticola'
$ws.Range("H8").Value2 = '(Schwein.) Murrill'
$ws.Range("Q8").Value2 = 331847
$ws.Range("R8").Value2 = 6626784
$ws.Range("L8").ClearContents()
$ws.Range("AC8").ClearContents()

# Row 9: A=111742070 (Lunglav)
$ws.Range("A9").Value2 = 111742070
$ws.Range("B9").Value2 = 78578
$ws.Range("D9").Value2 = 'NT'
$ws.Range("E9").Value2 = 6458
$ws.Range("F9").Value2 = 'Lunglav'
$ws.Range("G9").Value2 = 'Lobaria pulmonaria'
$ws.Range("H9").Value2 = '(L.) Hoffm.'
$ws.Range("Q9").Value2 = 331735
$ws.Range("R9").Value2 = 6626821
$ws.Range("L9").ClearContents()
$ws.Range("AC9").Value2 = 'På rönn'

# Row 10: A=111742096 (Skogstrappmossa)
$ws.Range("A10").Value2 = 111742096
$ws.Range("B10").Value2 = 94125
$ws.Range("D10").Value2 = 'NT'
$ws.Range("E10").Value2 = 54
$ws.Range("F10").Value2 = 'Skogstrappmossa'
$ws.Range("G10").Value2 = 'Anastrophyllum michauxii'
$ws.Range("H10").Value2 = '(F.Weber.) H.Buch'
$ws.Range("Q10").Value2 = 331780
$ws.Range("R10").Value2 = 6626798
$ws.Range("L10").Value2 = "'"
$ws.Range("AC10").Value2 = 'På både ved och på lodyta'

# Row 11: A=111742151 (Lopplummer)
$ws.Range("A11").Value2 = 111742151
$ws.Range("B11").Value2 = 95524
$ws.Range("D11").Value2 = 'LC'
$ws.Range("E11").Value2 = 221944
$ws.Range("F11").Value2 = 'Lopplummer'
$ws.Range("G11").Value2 = 'Huperzia selago'
$ws.Range("H11").Value2 = '(L.) Bernh. ex Schrank & Mart.'
$ws.Range("Q11").Value2 = 331815
$ws.Range("R11").Value2 = 6626779
$ws.Range("L11").Value2 = "'"
$ws.Range("AC11").ClearContents()

# Row 12: A=111742184 (Grov fjädermossa)
$ws.Range("A12").Value2 = 111742184
$ws.Range("B12").Value2 = 93159
$ws.Range("D12").Value2 = 'LC'
$ws.Range("E12").Value2 = 2666
$ws.Range("F12").Value2 = 'Grov fjädermossa'
$ws.Range("G12").Value2 = 'Neckera crispa'
$ws.Range("H12").Value2 = 'Hedw.'
$ws.Range("Q12").Value2 = 331834
$ws.Range("R12").Value2 = 6626785
$ws.Range("L12").Value2 = "'"
$ws.Range("AC12").Value2 = 'I bergsbrant'

# Row 13: A=111742077 (Stuplav)
$ws.Range("A13").Value2 = 111742077
$ws.Range("B13").Value2 = 78605
$ws.Range("D13").Value2 = 'LC'
$ws.Range("E13").Value2 = 6462
$ws.Range("F13").Value2 = 'Stuplav'
$ws.Range("G13").Value2 = 'Nephroma bellum'
$ws.Range("H13").Value2 = '(Spreng.) Tuck.'
$ws.Range("Q13").Value2 = 331735
$ws.Range("R13").Value2 = 6626821
$ws.Range("L13").ClearContents()
$ws.Range("AC13").Value2 = 'På rönn'

# Row 14: A=111742181 (Platt fjädermossa)
$ws.Range("A14").Value2 = 111742181
$ws.Range("B14").Value2 = 93158
$ws.Range("D14").Value2 = 'LC'
$ws.Range("E14").Value2 = 2667
$ws.Range("F14").Value2 = 'Platt fjädermossa'
$ws.Range("G14").Value2 = 'Neckera complanata'
$ws.Range("H14").Value2 = '(Hedw.) Huebener'
$ws.Range("Q14").Value2 = 331834
$ws.Range("R14").Value2 = 6626785
$ws.Range("L14").Value2 = "'"
$ws.Range("AC14").Value2 = 'I bergsbrant'
